$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply AutoFilter over the current (pre-edit) data range A1:F8 so that
# subsequently inserted rows below it are not swept into the filter range.
$ws.AutoFilterMode = $false
$ws.Range("A1:F8").AutoFilter() | Out-Null

# Duplicate row 8 twice (as rows 9 and 10), using Copy + Insert so that the
# row formatting (including the "empty" alignment/fill style used by column E)
# is carried over faithfully.
$ws.Range("A8:F8").Copy()
$ws.Rows(9).Insert()

$ws.Range("A8:F8").Copy()
$ws.Rows(10).Insert()

# Fill in the new test-case text for the two appended rows.
$ws.Range("A9").Value = "Product_Summary-Edit_Product_view-edit_name_of_account_[WEB]_1"
$ws.Range("B9").Value = "C70767"

$ws.Range("A10").Value = "Product_Summary-Edit_Product_view-edit_name_of_account-Invalid_[WEB]_1"
$ws.Range("B10").Value = "C70768"

# Update the selection to match the author's final cursor position.
$ws.Range("A10").Select() | Out-Null

# Update the hidden _FilterDatabase defined name to match the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$8"
    }
}

# Clear the stale row outline-level summary while keeping the column one.
$ws.Outline.ShowLevels(0)
